# modified born position of city
# Update the RelivePos (column E) for the villageScene row (row 2)
# from "0,0,0" to "20,0,-137".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "20,0,-137"
